# AI_Timesheets.xlsx update: new employee (Ray Beer) and the week's dates
# shifted back three weeks (the previous data pointed at the wrong week).
#
# Dates are entered through a literal text formula and then pasted back as
# values so Excel stores them as plain text (matching the workbook's
# existing "MM/DD/YYYY" text entries) instead of auto-converting them into
# date serial numbers, while leaving each cell's existing formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($range, $text)
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Employee name (D6)
Set-TextValue $ws.Range("D6") "Ray Beer"

# Date column (C10:C16) - same week-day pattern, three weeks earlier
Set-TextValue $ws.Range("C10") "10/28/2024"
Set-TextValue $ws.Range("C11") "10/29/2024"
Set-TextValue $ws.Range("C12") "10/30/2024"
Set-TextValue $ws.Range("C13") "10/31/2024"
Set-TextValue $ws.Range("C14") "11/01/2024"
Set-TextValue $ws.Range("C15") "11/02/2024"
Set-TextValue $ws.Range("C16") "11/03/2024"

$null = $ws.Range("A1").Select()
